$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data: new record MCH199 ---
$ws.Range("A2").Value = "MCH199"
$ws.Range("C2").Value = "PACKAGE ON HUMAN RIGHTS TO COINCIDE WITH THE VISIT OF ARCHBISHOP DESMOND TUTU TO BIRMINGHAM 1988"

# date_s must be stored as text (not a number) - force Text format before entry
$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "1988"

$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# --- Formatting: row 2 gets a dedicated Calibri 10pt style ---
$row2 = $ws.Range("A2:H2")
$row2.Font.Name = "Calibri"
$row2.Font.Size = 10
$row2.Font.ThemeColor = 1

# Re-sync D2's number format back to General (inherits from A2, which already
# carries the final row-2 font) now that its value is safely stored as text.
$ws.Range("A2").Copy()
$d2.PasteSpecial(-4122)

# B2 (alternativeIdentifiers) has no value for this record
$ws.Range("B2").Clear()

# --- Sheet view: header row frozen, A2:I2 selected ---
$ws.Range("A2:I2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
